$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56
$ws.Range("A56").Value = 111840155
$ws.Range("B56").Value = 78578
$ws.Range('D56').Value = 'NT'
$ws.Range("E56").Value = 6458
$ws.Range('F56').Value = 'Lunglav'
$ws.Range('G56').Value = 'Lobaria pulmonaria'
$ws.Range('H56').Value = '(L.) Hoffm.'
$ws.Range("I56").ClearContents()
$ws.Range("Q56").Value = 611724.3624407195
$ws.Range("R56").Value = 7036587.505142178

# Row 57
$ws.Range("A57").Value = 111840166
$ws.Range("B57").Value = 96348
$ws.Range('D57').Value = 'VU'
$ws.Range("E57").Value = 220787
$ws.Range('F57').Value = 'Knärot'
$ws.Range('G57').Value = 'Goodyera repens'
$ws.Range('H57').Value = '(L.) R. Br.'
$ws.Range('I57').Value = '5'
$ws.Range("Q57").Value = 611872.0517480521
$ws.Range("R57").Value = 7036366.032528495

# Row 58
$ws.Range("A58").Value = 111840113
$ws.Range("B58").Value = 89405
$ws.Range('D58').Value = 'NT'
$ws.Range("E58").Value = 1202
$ws.Range('F58').Value = 'Ullticka'
$ws.Range('G58').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H58').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I58").ClearContents()
$ws.Range("Q58").Value = 611968.6292715519
$ws.Range("R58").Value = 7036400.323734847

# Row 59
$ws.Range("A59").Value = 111840117
$ws.Range("B59").Value = 89405
$ws.Range('D59').Value = 'NT'
$ws.Range("E59").Value = 1202
$ws.Range('F59').Value = 'Ullticka'
$ws.Range('G59').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H59').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I59").ClearContents()
$ws.Range("Q59").Value = 611725.3738043295
$ws.Range("R59").Value = 7036571.416293882

# Row 60
$ws.Range("A60").Value = 111840114
$ws.Range("B60").Value = 89405
$ws.Range('D60').Value = 'NT'
$ws.Range("E60").Value = 1202
$ws.Range('F60').Value = 'Ullticka'
$ws.Range('G60').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H60').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I60").ClearContents()
$ws.Range("Q60").Value = 611931.1515982752
$ws.Range("R60").Value = 7036433.051986331

# Row 61
$ws.Range("A61").Value = 111840141
$ws.Range("B61").Value = 78605
$ws.Range('D61').Value = 'LC'
$ws.Range("E61").Value = 6462
$ws.Range('F61').Value = 'Stuplav'
$ws.Range('G61').Value = 'Nephroma bellum'
$ws.Range('H61').Value = '(Spreng.) Tuck.'
$ws.Range("I61").ClearContents()
$ws.Range("Q61").Value = 611995.7935480368
$ws.Range("R61").Value = 7036214.038866865

# Row 62
$ws.Range("A62").Value = 111840164
$ws.Range("B62").Value = 96348
$ws.Range('D62').Value = 'VU'
$ws.Range("E62").Value = 220787
$ws.Range('F62').Value = 'Knärot'
$ws.Range('G62').Value = 'Goodyera repens'
$ws.Range('H62').Value = '(L.) R. Br.'
$ws.Range('I62').Value = '50'
$ws.Range("Q62").Value = 611550.5735735258
$ws.Range("R62").Value = 7036580.982027167

# Row 63
$ws.Range("A63").Value = 111840170
$ws.Range("B63").Value = 96348
$ws.Range('D63').Value = 'VU'
$ws.Range("E63").Value = 220787
$ws.Range('F63').Value = 'Knärot'
$ws.Range('G63').Value = 'Goodyera repens'
$ws.Range('H63').Value = '(L.) R. Br.'
$ws.Range('I63').Value = '20'
$ws.Range("Q63").Value = 611721.2521968643
$ws.Range("R63").Value = 7036599.489451895
